$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - add new columns I (I0) and J (IF)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold, border, centered) from H1 onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data values for columns I (I0) and J (IF), rows 2-19
$data = @(
    @(7, 8),
    @(9, 9),
    @(8, 8),
    @(7, 8),
    @(8, 8),
    @(6, 6),
    @(6, 6),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(5, 5),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(7, 7),
    @(6, 6),
    @(4, 4),
    @(7, 7)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
